$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Usuario column (A) to raul.gomez for all data rows
$ws.Range("A2").Value = "raul.gomez"
$ws.Range("A3").Value = "raul.gomez"
$ws.Range("A4").Value = "raul.gomez"
$ws.Range("A5").Value = "raul.gomez"

# Update event names in column B
$ws.Range("B2").Value = "Presentacion 28"
$ws.Range("B3").Value = "Junta con el gato"
$ws.Range("B4").Value = "Curso Scrum"
$ws.Range("B5").Value = "Curso excel"

# Update date in C2
$ws.Range("C2").Value = (Get-Date -Year 2018 -Month 2 -Day 26 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)

# Move the special (2nd) cell style from A4/A5 down to D4/D5:
# copy format only from A4/A5 onto D4/D5, then reset A4/A5 to the default style.
$ws.Range("A4").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("A5").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A4").Style = "Normal"
$ws.Range("A5").Style = "Normal"

$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).AutoFit()

# Update selection to D5
$ws.Range("D5").Select()
